$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.941.64"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.480.91"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.22"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.58"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "4.030.85"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.48"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "3.501.72"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "62.897.86"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("E25").Value = "  +12.33%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.73"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.50%  "
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +8.05%  "
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "0.0₃0587"
$ws.Range("E47").Value = "  +34.23%  "
$ws.Range("E48").Value = "  +11.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.39"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.16"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  -3.69%  "
